# Apply the cryptos-list price/volume refresh described by the commit diff.
# Values are written as TEXT (matching the source workbook, which stores every
# D/E cell as an inline string) by forcing a Text number format for the duration
# of the write, then clearing formats again so the cell keeps style index 0 -
# exactly like the original cells - instead of picking up a "@" text style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$edits = @(
    @{ Cell = 'D2'; Value = '94.481.95' },
    @{ Cell = 'E2'; Value = '  -1.45%  ' },
    @{ Cell = 'D3'; Value = '3.333.37' },
    @{ Cell = 'E3'; Value = '  -3.92%  ' },
    @{ Cell = 'D4'; Value = '1.00' },
    @{ Cell = 'E4'; Value = '  -0.14%  ' },
    @{ Cell = 'D5'; Value = '231.93' },
    @{ Cell = 'E5'; Value = '  -4.17%  ' },
    @{ Cell = 'D6'; Value = '620.41' },
    @{ Cell = 'E6'; Value = '  -3.95%  ' },
    @{ Cell = 'D7'; Value = '1.39' },
    @{ Cell = 'E7'; Value = '  -5.06%  ' },
    @{ Cell = 'E8'; Value = '  -5.89%  ' },
    @{ Cell = 'E9'; Value = '  -0.10%  ' },
    @{ Cell = 'D10'; Value = '0.942' },
    @{ Cell = 'E10'; Value = '  -5.29%  ' },
    @{ Cell = 'D11'; Value = '3.334.12' },
    @{ Cell = 'E11'; Value = '  -3.90%  ' },
    @{ Cell = 'D12'; Value = '42.22' },
    @{ Cell = 'E12'; Value = '  -1.21%  ' },
    @{ Cell = 'E13'; Value = '  -2.87%  ' },
    @{ Cell = 'B14'; Value = 'WrappedBTC' },
    @{ Cell = 'C14'; Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc' },
    @{ Cell = 'D14'; Value = '94.224.01' },
    @{ Cell = 'E14'; Value = '  -1.53%  ' },
    @{ Cell = 'B15'; Value = 'Toncoin' },
    @{ Cell = 'C15'; Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton' },
    @{ Cell = 'D15'; Value = '5.99' },
    @{ Cell = 'E15'; Value = '  -2.44%  ' },
    @{ Cell = 'D16'; Value = '3.951.00' },
    @{ Cell = 'E16'; Value = '  -4.24%  ' },
    @{ Cell = 'D17'; Value = '0.0000245' },
    @{ Cell = 'E17'; Value = '  -3.94%  ' },
    @{ Cell = 'D18'; Value = '8.14' },
    @{ Cell = 'E18'; Value = '  -3.76%  ' },
    @{ Cell = 'D19'; Value = '3.333.67' },
    @{ Cell = 'E19'; Value = '  -3.92%  ' },
    @{ Cell = 'D20'; Value = '17.30' },
    @{ Cell = 'E20'; Value = '  -5.28%  ' },
    @{ Cell = 'E21'; Value = '  -7.09%  ' },
    @{ Cell = 'D22'; Value = '3.52' },
    @{ Cell = 'E22'; Value = '  +8.57%  ' },
    @{ Cell = 'D23'; Value = '496.22' },
    @{ Cell = 'E23'; Value = '  -2.73%  ' },
    @{ Cell = 'E24'; Value = '  -11.15%  ' },
    @{ Cell = 'D25'; Value = '0.0000183' },
    @{ Cell = 'E25'; Value = '  -5.18%  ' },
    @{ Cell = 'D26'; Value = '6.22' },
    @{ Cell = 'E26'; Value = '  -6.44%  ' },
    @{ Cell = 'D27'; Value = '90.04' },
    @{ Cell = 'E27'; Value = '  -2.41%  ' },
    @{ Cell = 'D28'; Value = '11.78' },
    @{ Cell = 'E28'; Value = '  -4.34%  ' },
    @{ Cell = 'D29'; Value = '3.508.15' },
    @{ Cell = 'E29'; Value = '  -4.19%  ' },
    @{ Cell = 'E30'; Value = '  +0.01%  ' },
    @{ Cell = 'D31'; Value = '11.21' },
    @{ Cell = 'E31'; Value = '  -5.30%  ' },
    @{ Cell = 'D32'; Value = '0.138' },
    @{ Cell = 'E32'; Value = '  +0.19%  ' },
    @{ Cell = 'E33'; Value = '  -4.27%  ' },
    @{ Cell = 'D34'; Value = '1.01' },
    @{ Cell = 'E34'; Value = '  +0.53%  ' },
    @{ Cell = 'D35'; Value = '0.175' },
    @{ Cell = 'E35'; Value = '  -5.15%  ' },
    @{ Cell = 'D36'; Value = '28.46' },
    @{ Cell = 'E36'; Value = '  -8.20%  ' },
    @{ Cell = 'D37'; Value = '0.533' },
    @{ Cell = 'E37'; Value = '  -7.28%  ' },
    @{ Cell = 'D38'; Value = '533.35' },
    @{ Cell = 'E38'; Value = '  +2.99%  ' },
    @{ Cell = 'B39'; Value = 'RenderToken' },
    @{ Cell = 'C39'; Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render' },
    @{ Cell = 'D39'; Value = '7.39' },
    @{ Cell = 'E39'; Value = '  -5.46%  ' },
    @{ Cell = 'B40'; Value = 'USDe' },
    @{ Cell = 'C40'; Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde' },
    @{ Cell = 'D40'; Value = '1.00' },
    @{ Cell = 'E40'; Value = '  +0.04%  ' },
    @{ Cell = 'D41'; Value = '1.37' },
    @{ Cell = 'E41'; Value = '  -5.79%  ' },
    @{ Cell = 'D42'; Value = '0.148' },
    @{ Cell = 'E42'; Value = '  -1.68%  ' },
    @{ Cell = 'D43'; Value = '0.873' },
    @{ Cell = 'E43'; Value = '  -4.53%  ' },
    @{ Cell = 'B44'; Value = 'MantraDAO' },
    @{ Cell = 'C44'; Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om' },
    @{ Cell = 'D44'; Value = '3.79' },
    @{ Cell = 'E44'; Value = '  +4.52%  ' },
    @{ Cell = 'B45'; Value = 'WhiteBITCoin' },
    @{ Cell = 'C45'; Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt' },
    @{ Cell = 'D45'; Value = '24.13' },
    @{ Cell = 'E45'; Value = '  -0.04%  ' },
    @{ Cell = 'D46'; Value = '0.0418' },
    @{ Cell = 'E46'; Value = '  -0.10%  ' },
    @{ Cell = 'D47'; Value = '1.69' },
    @{ Cell = 'E47'; Value = '  -1.33%  ' },
    @{ Cell = 'D48'; Value = '5.43' },
    @{ Cell = 'E48'; Value = '  -2.60%  ' },
    @{ Cell = 'D49'; Value = '53.41' },
    @{ Cell = 'E49'; Value = '  -0.38%  ' },
    @{ Cell = 'D50'; Value = '2.11' },
    @{ Cell = 'E50'; Value = '  -3.10%  ' },
    @{ Cell = 'D51'; Value = '8.05' },
    @{ Cell = 'E51'; Value = '  -1.59%  ' }
)

foreach ($edit in $edits) {
    $c = $ws.Range($edit.Cell)
    $c.NumberFormat = "@"
    $c.Value = $edit.Value
    $c.ClearFormats()
}
